{"js": "// Update the worksheet date and all 25 division problems/answers to the\n// new day's content. Every line of text in the document is unique, both\n// before and after the edit, so a sequence of exact search+replace calls\n// (processed in this fixed order) deterministically reaches the target\n// state even though one new value (\"26\u00f75=5, 1\") happens to equal an old\n// value that an earlier step has already replaced away.\nconst replacements = [\n  [\"2024-09-20 Friday\", \"2024-09-21 Saturday\"],\n  [\"77\u00f76=12, 5\", \"40\u00f74=10, 0\"],\n  [\"35\u00f78=4, 3\", \"45\u00f73=15, 0\"],\n  [\"26\u00f75=5, 1\", \"15\u00f74=3, 3\"],\n  [\"83\u00f76=13, 5\", \"10\u00f72=5, 0\"],\n  [\"40\u00f78=5, 0\", \"78\u00f78=9, 6\"],\n  [\"39\u00f73=13, 0\", \"25\u00f78=3, 1\"],\n  [\"85\u00f75=17, 0\", \"53\u00f72=26, 1\"],\n  [\"76\u00f77=10, 6\", \"46\u00f79=5, 1\"],\n  [\"13\u00f73=4, 1\", \"36\u00f74=9, 0\"],\n  [\"17\u00f74=4, 1\", \"68\u00f73=22, 2\"],\n  [\"15\u00f79=1, 6\", \"19\u00f76=3, 1\"],\n  [\"51\u00f78=6, 3\", \"43\u00f79=4, 7\"],\n  [\"94\u00f79=10, 4\", \"26\u00f75=5, 1\"],\n  [\"61\u00f79=6, 7\", \"79\u00f72=39, 1\"],\n  [\"70\u00f77=10, 0\", \"14\u00f77=2, 0\"],\n  [\"18\u00f72=9, 0\", \"57\u00f78=7, 1\"],\n  [\"90\u00f77=12, 6\", \"92\u00f74=23, 0\"],\n  [\"75\u00f73=25, 0\", \"83\u00f79=9, 2\"],\n  [\"92\u00f77=13, 1\", \"46\u00f72=23, 0\"],\n  [\"71\u00f72=35, 1\", \"18\u00f76=3, 0\"],\n  [\"82\u00f77=11, 5\", \"36\u00f73=12, 0\"],\n  [\"55\u00f79=6, 1\", \"55\u00f76=9, 1\"],\n  [\"28\u00f79=3, 1\", \"98\u00f79=10, 8\"],\n  [\"62\u00f75=12, 2\", \"26\u00f79=2, 8\"],\n  [\"62\u00f78=7, 6\", \"90\u00f74=22, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 division problems/answers to the\n# new day's content. Every line of text in the document is unique, both\n# before and after the edit, so a sequence of exact Find/Replace calls\n# (processed in this fixed order) deterministically reaches the target\n# state even though one new value (\"26\u00f75=5, 1\") happens to equal an old\n# value that an earlier step has already replaced away.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-09-20 Friday\", \"2024-09-21 Saturday\"),\n    @(\"77\u00f76=12, 5\", \"40\u00f74=10, 0\"),\n    @(\"35\u00f78=4, 3\", \"45\u00f73=15, 0\"),\n    @(\"26\u00f75=5, 1\", \"15\u00f74=3, 3\"),\n    @(\"83\u00f76=13, 5\", \"10\u00f72=5, 0\"),\n    @(\"40\u00f78=5, 0\", \"78\u00f78=9, 6\"),\n    @(\"39\u00f73=13, 0\", \"25\u00f78=3, 1\"),\n    @(\"85\u00f75=17, 0\", \"53\u00f72=26, 1\"),\n    @(\"76\u00f77=10, 6\", \"46\u00f79=5, 1\"),\n    @(\"13\u00f73=4, 1\", \"36\u00f74=9, 0\"),\n    @(\"17\u00f74=4, 1\", \"68\u00f73=22, 2\"),\n    @(\"15\u00f79=1, 6\", \"19\u00f76=3, 1\"),\n    @(\"51\u00f78=6, 3\", \"43\u00f79=4, 7\"),\n    @(\"94\u00f79=10, 4\", \"26\u00f75=5, 1\"),\n    @(\"61\u00f79=6, 7\", \"79\u00f72=39, 1\"),\n    @(\"70\u00f77=10, 0\", \"14\u00f77=2, 0\"),\n    @(\"18\u00f72=9, 0\", \"57\u00f78=7, 1\"),\n    @(\"90\u00f77=12, 6\", \"92\u00f74=23, 0\"),\n    @(\"75\u00f73=25, 0\", \"83\u00f79=9, 2\"),\n    @(\"92\u00f77=13, 1\", \"46\u00f72=23, 0\"),\n    @(\"71\u00f72=35, 1\", \"18\u00f76=3, 0\"),\n    @(\"82\u00f77=11, 5\", \"36\u00f73=12, 0\"),\n    @(\"55\u00f79=6, 1\", \"55\u00f76=9, 1\"),\n    @(\"28\u00f79=3, 1\", \"98\u00f79=10, 8\"),\n    @(\"62\u00f75=12, 2\", \"26\u00f79=2, 8\"),\n    @(\"62\u00f78=7, 6\", \"90\u00f74=22, 2\")\n)\n\nforeach ($pair in $pairs) {\n    $rng = $d.Content\n    $rng.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}\n"}
